$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '28.135.62'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +2.35%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.911.19'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +2.08%  '

$ws.Range('E4').Value = '  -1.25%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '316.39'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.15%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.004'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -1.05%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4833'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +1.11%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3819'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +1.25%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07361'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -0.14%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.9362'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -0.15%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '20.82'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.47%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.07792'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -0.48%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.903.34'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +0.88%  '

$ws.Range('E14').Value = '  +1.37%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '6.636'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +0.81%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '91.48'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +0.67%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '1.006'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -1.23%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.000008840'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -0.65%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '1.004'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -1.00%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '28.160.83'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +2.28%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '14.87'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.31%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.155'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.47%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '2.148.41'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +1.77%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '10.91'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +1.77%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '156.70'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +1.69%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.914'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -2.40%  '

$ws.Range('E27').Value = '  +0.21%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.115'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +4.86%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '116.33'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +0.36%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '4.952'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -0.76%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.08926'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -0.11%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.351'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -0.01%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.255'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +2.99%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.7701'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +2.48%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '4.685'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +1.60%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.612'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -2.86%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.02051'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -0.14%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.103'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -1.27%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.05316'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.26%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.5508'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +3.01%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.979'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.95%  '

$ws.Range('E43').Value = '  -0.17%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '8.469'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +0.59%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '10.72'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.57%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.4840'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +0.34%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '107.50'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +4.32%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.005'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -1.05%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.660'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.10%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '68.45'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +1.88%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.06108'
$ws.Range('D51').Style = 'Normal'
